# DMZ fixed to come from LAN switch rather than firewall.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Address Chart")

# Row 6 previously described the DMZ hand-off as coming from the PA Firewall
# (Device=PA Firewall, Interface=E1/3, Gateway=N/A, Zone=DMZ). Correct it so
# the DMZ hand-off instead comes from the LAN Switch via VLAN100.
$ws.Range("B6").Value = "LAN Switch"
$ws.Range("C6").Value = "VLAN100"
$ws.Range("F6").Value = "192.168.150.1"
$ws.Range("G6").Value = "DMZ "

# Note on the DMZ Server row that it is reached via VLAN100.
$ws.Range("H8").Value = "VLAN100"

# Keep the previously selected cell as the active selection.
$ws.Range("G5").Select()
